# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Pepino ensalada" above the existing row 308, pushing all subsequent rows
# (previously 308-383) down by one (to 309-384).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 308..383 down to 309..384, leaving a blank row 308 (formatting
# of the row above carries down automatically, same as Excel's native
# Insert behaviour).
$ws.Rows.Item(308).Insert()

# Populate the newly-inserted row 308 with the new weekly observation.
$ws.Cells.Item(308, 1).Value = 4
$ws.Cells.Item(308, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(308, 3).Value = "Los Lagos"
$ws.Cells.Item(308, 4).Value = 44964
$ws.Cells.Item(308, 5).Value = 10
$ws.Cells.Item(308, 6).Value = 100112043
$ws.Cells.Item(308, 7).Value = "Pepino ensalada"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 400
$ws.Cells.Item(308, 11).Value = 13000
$ws.Cells.Item(308, 12).Value = 14000
$ws.Cells.Item(308, 13).Value = 13500
$ws.Cells.Item(308, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(308, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(308, 16).Value = 225
$ws.Cells.Item(308, 17).Value = 60
$ws.Cells.Item(308, 18).Value = "Hortaliza"
